$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "time_taken" column, using the same style as the
# existing header row (bold, bordered, centered - reuses style index 1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Per-row timestamps recorded for the "time_taken" metadata column.
$timestamps = @(
    "2021-10-05 10:51:04.890415",
    "2021-10-05 10:51:04.890428",
    "2021-10-05 10:51:04.890432",
    "2021-10-05 10:51:04.890435",
    "2021-10-05 10:51:04.890438",
    "2021-10-05 10:51:04.890442",
    "2021-10-05 10:51:04.890445",
    "2021-10-05 10:51:04.890448",
    "2021-10-05 10:51:04.890451",
    "2021-10-05 10:51:04.890454",
    "2021-10-05 10:51:04.890457"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
